# Weekly crime-data refresh: new report week (11/7/2022-11/13/2022, Vol 29 No. 45)
# and updated CompStat figures for rows 15-29.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (partial run replace via Characters, 1-based) ---
# A8 shared string: "Volume 29   Number  44" -> "...45" (chars 21-22 = "44")
$ws.Range("A8").Characters(21, 2).Text = "45"

# C9 shared string: "Report Covering the Week  10/31/2022  Through  11/6/2022"
# First replace shifts the string, so the second offset is computed post-edit.
$ws.Range("C9").Characters(27, 10).Text = "11/7/2022"
$ws.Range("C9").Characters(47, 9).Text = "11/13/2022"

# --- Crime statistics table updates (rows 15-29) ---

# Row 15
$ws.Range("C15").Value = 1.0
$ws.Range("C15").NumberFormat = "#,##0"
$ws.Range("F15").Value = 1.0
$ws.Range("F15").NumberFormat = "#,##0"
$ws.Range("I15").Value = 19.0
$ws.Range("K15").Value = 58.333333333333
$ws.Range("L15").Value = 58.333333333333
$ws.Range("M15").Value = 11.764705882352
$ws.Range("N15").Value = -5.0

# Row 16
$ws.Range("C16").Value = 6.0
$ws.Range("D16").Value = 2.0
$ws.Range("E16").Value = 200.0
$ws.Range("F16").Value = 16.0
$ws.Range("H16").Value = 33.333333333333
$ws.Range("I16").Value = 178.0
$ws.Range("J16").Value = 116.0
$ws.Range("K16").Value = 53.448275862069
$ws.Range("L16").Value = 9.876543209876
$ws.Range("M16").Value = -22.943722943722
$ws.Range("N16").Value = -79.326364692218

# Row 17
$ws.Range("C17").Value = 5.0
$ws.Range("D17").Value = 6.0
$ws.Range("E17").Value = -16.666666666666
$ws.Range("G17").Value = 25.0
$ws.Range("H17").Value = 16.0
$ws.Range("I17").Value = 266.0
$ws.Range("J17").Value = 226.0
$ws.Range("K17").Value = 17.699115044247
$ws.Range("L17").Value = 2.702702702702
$ws.Range("M17").Value = 37.823834196891
$ws.Range("N17").Value = -0.374531835205

# Row 18
$ws.Range("C18").Value = 3.0
$ws.Range("D18").Value = 10.0
$ws.Range("E18").Value = -70.0
$ws.Range("F18").Value = 14.0
$ws.Range("G18").Value = 27.0
$ws.Range("H18").Value = -48.148148148148
$ws.Range("I18").Value = 236.0
$ws.Range("J18").Value = 198.0
$ws.Range("K18").Value = 19.191919191919
$ws.Range("L18").Value = -15.107913669064
$ws.Range("M18").Value = -40.253164556962
$ws.Range("N18").Value = -86.483390607102

# Row 19
$ws.Range("C19").Value = 10.0
$ws.Range("D19").Value = 10.0
$ws.Range("E19").Value = 0.0
$ws.Range("F19").Value = 57.0
$ws.Range("G19").Value = 46.0
$ws.Range("H19").Value = 23.91304347826
$ws.Range("I19").Value = 566.0
$ws.Range("J19").Value = 467.0
$ws.Range("K19").Value = 21.19914346895
$ws.Range("L19").Value = 8.221797323135
$ws.Range("M19").Value = 50.933333333333
$ws.Range("N19").Value = -0.176366843033

# Row 20
$ws.Range("C20").Value = 10.0
$ws.Range("D20").Value = 8.0
$ws.Range("E20").Value = 25.0
$ws.Range("F20").Value = 27.0
$ws.Range("G20").Value = 22.0
$ws.Range("H20").Value = 22.727272727272
$ws.Range("I20").Value = 269.0
$ws.Range("J20").Value = 187.0
$ws.Range("K20").Value = 43.850267379679
$ws.Range("L20").Value = 73.548387096774
$ws.Range("M20").Value = -16.19937694704
$ws.Range("N20").Value = -91.237785016286

# Row 21
$ws.Range("C21").Value = 35.0
$ws.Range("D21").Value = 36.0
$ws.Range("E21").Value = -2.777777777777
$ws.Range("F21").Value = 144.0
$ws.Range("G21").Value = 132.0
$ws.Range("H21").Value = 9.090909090909
$ws.Range("I21").Value = 1536.0
$ws.Range("J21").Value = 1208.0
$ws.Range("K21").Value = 27.152317880794
$ws.Range("L21").Value = 10.583153347732
$ws.Range("M21").Value = 0.0
$ws.Range("N21").Value = -76.542455711667

# Row 22
$ws.Range("L22").Value = -31.578947368421

# Row 24
$ws.Range("C24").Value = 36.0
$ws.Range("D24").Value = 28.0
$ws.Range("E24").Value = 28.571428571428
$ws.Range("F24").Value = 126.0
$ws.Range("G24").Value = 108.0
$ws.Range("H24").Value = 16.666666666666
$ws.Range("I24").Value = 1279.0
$ws.Range("J24").Value = 1165.0
$ws.Range("K24").Value = 9.785407725321
$ws.Range("L24").Value = -0.15612802498
$ws.Range("M24").Value = 39.021739130434

# Row 25
$ws.Range("C25").Value = 5.0
$ws.Range("D25").Value = 11.0
$ws.Range("E25").Value = -54.545454545454
$ws.Range("F25").Value = 42.0
$ws.Range("G25").Value = 38.0
$ws.Range("H25").Value = 10.526315789473
$ws.Range("I25").Value = 470.0
$ws.Range("J25").Value = 439.0
$ws.Range("K25").Value = 7.061503416856
$ws.Range("L25").Value = 17.5
$ws.Range("M25").Value = -24.920127795527

# Row 26
$ws.Range("F26").Value = 4.0
$ws.Range("I26").Value = 31.0
$ws.Range("K26").Value = 72.222222222222
$ws.Range("L26").Value = 24.0

# Row 27
$ws.Range("D27").Value = 2.0
$ws.Range("D27").NumberFormat = "#,##0"
$ws.Range("E27").Value = -100.0
$ws.Range("E27").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("F27").Value = 3.0
$ws.Range("G27").Value = 3.0
$ws.Range("H27").Value = 0.0
$ws.Range("J27").Value = 47.0
$ws.Range("K27").Value = 34.042553191489

# Row 28
$ws.Range("N28").Value = -69.230769230769

# Row 29
$ws.Range("N29").Value = -69.565217391304
